$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of GILD "Random" prediction data to append (rows 7-18)
$rows = @(
    @(42600.782094907408, 13, 87, 70, 30),
    @(42600.804745370369, 57, 43, 52, 48),
    @(42600.806909722225, 99, 1, 67, 33),
    @(42600.823356481480, 37, 63, 53, 47),
    @(42600.830381944441, 39, 61, 10, 90),
    @(42600.841898148145, 42, 58, 48, 52),
    @(42600.861435185187, 95, 5, 46, 54),
    @(42600.868553240740, 41, 59, 87, 13),
    @(42600.878773148150, 68, 32, 45, 55),
    @(42600.884097222224, 24, 76, 50, 50),
    @(42600.885451388887, 45, 55, 69, 31),
    @(42600.886643518519, 63, 37, 83, 17)
)

# Carry the date-formatted style from the last existing row (A6) down through
# the newly appended rows (A7:A18) before writing the values.
$ws.Range("A6").Copy($ws.Range("A7:A18"))

$r = 7
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]

    $ws.Cells.Item($r, 2).Value = "Random"

    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 5).Value = 0
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 7).Value = 0

    $ws.Cells.Item($r, 8).Value = $row[1]
    $ws.Cells.Item($r, 9).Value = $row[2]

    $ws.Cells.Item($r, 10).Value = 0
    $ws.Cells.Item($r, 11).Value = 0

    $ws.Cells.Item($r, 12).Value = $row[3]
    $ws.Cells.Item($r, 13).Value = $row[4]

    $r++
}

# Widen column A to fit the additional data
$ws.Columns.Item(1).ColumnWidth = 14
